# Reformat the death event parameters in the "parameter_values" sheet so
# that the linear model can be used to calculate an individual's risk of
# death (treatment effect is incorporated as a predictor of death).
#
# Concretely:
#  1. Remove the "cfr_obstructed_labour" row (old row 35) - this value is
#     no longer used now that the linear model handles this risk, and
#     everything below shifts up by one row.
#  2. Replace the two rows "prob_deliver_ventouse" and
#     "prob_deliver_forceps" with a single new row
#     "prob_successful_assisted_vaginal_delivery" (value 0.7), which is
#     used as a predictor in the new linear model.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Delete the cfr_obstructed_labour row ---------------------------
$ws.Rows.Item(35).Delete()

# After the deletion above, the rows that used to be 80/81
# ("prob_deliver_ventouse" / "prob_deliver_forceps") are now 79/80.

# --- 2. Collapse the two "prob_deliver_*" rows into one -----------------
$ws.Rows.Item(80).Delete()

$ws.Range("A79").ClearFormats()
$ws.Range("A79").Value = "prob_successful_assisted_vaginal_delivery"
$ws.Range("B79").Value = 0.7

# --- 3. Update the view so the newly edited area is visible -------------
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
$ws.Rows.Item(35).Select()
